# chore: adapt column header formatting to respective input file names (#7)
#
# Renames the "_old"/"_new" column header suffixes to the concrete AHB
# format-version suffixes ("_FV2304" / "_FV2310"), wraps the data range in a
# real Excel Table (so the header renames / filtering survive round-trips),
# and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (A1:U1): "_old" -> "_FV2304", "_new" -> "_FV2310" ---
$oldSuffixHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$col = 1
foreach ($name in $oldSuffixHeaders) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2304"
    $col++
}

# Column K ("diff") is unchanged.
$col++

foreach ($name in $oldSuffixHeaders) {
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2310"
    $col++
}

# --- 2) Turn the used range into a real Table so the renamed headers ---
#        become the table's column names.
$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
